# "Generate Report for Handback" - refresh the handback status timestamps
# (and a priority value) across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-24 10:20:29"
$wsOverview.Range("G3").Value = "2016-10-24 10:20:29"

# --- zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-10-24 10:20:17"
$wsZhCn.Range("H3").Value = "2016-10-24 10:20:17"
$wsZhCn.Range("K2").Value = "2016-10-24 10:21:00"
$wsZhCn.Range("K3").Value = "2016-10-24 10:21:00"

# --- de-de sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-10-24 10:20:29"
$wsDeDe.Range("H3").Value = "2016-10-24 10:20:29"
$wsDeDe.Range("K2").Value = "2016-10-24 10:21:17"
$wsDeDe.Range("K3").Value = "2016-10-24 10:21:17"
